# "Add tutorials for condition monitoring application."
#
# The log sheet (Sheet1) gets one new test-run entry appended as row 8:
#   A8 = test id "20240321_122650"
#   B8 = description of the failed pick-up attempt
#
# Column B already auto-sizes to fit its longest value ("bestFit"); adding a
# longer description widens it further, so we resize it to fit the new text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "20240321_122650"
$ws.Range("B8").Value = "Try to move to pick up an item from second floor of a shelf. But failed."

# Re-fit column B now that it holds a longer description.
$ws.Columns.Item(2).ColumnWidth = 63.45

# Leave the selection where the author ended up after entering the data.
$ws.Range("B16").Select()
